# Applies the text edits described by the commit diff:
#  - Slide 1: title quote changes from "Team Bahug Taler" to "Itghurls"
#  - Slide 2: subtitle line changes to "PitchItup - an auto generated pitchdeck"
#  - Slide 3: problem bullets updated ("- Having hard time creating pitchdeck" / "- ")
#  - Slide 4: subtitle line changes to "PitchItup - an auto generated pitchdeck"
#
# Shapes in this deck use spAutoFit, so pushing new text through the COM
# TextRange re-flows (and resizes) the shape. We snapshot/restore each
# shape's Width/Height around the edit so only the run text itself changes.

function Set-ShapeRunText($shape, [int]$runIndex, [string]$newText) {
    $origWidth = $shape.Width
    $origHeight = $shape.Height
    $tr = $shape.TextFrame.TextRange
    $tr.Runs($runIndex, 1).Text = $newText
    $shape.Width = $origWidth
    $shape.Height = $origHeight
}

$p = $ppt.ActivePresentation

# --- Slide 1 ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$shape1 = $s1.Shapes.Item(2)
Set-ShapeRunText $shape1 1 '"Itghurls"'

# --- Slide 2 ---------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shape2 = $s2.Shapes.Item(2)
Set-ShapeRunText $shape2 2 "PitchItup - an auto generated pitchdeck"

# --- Slide 3 ---------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$shape3 = $s3.Shapes.Item(3)
$origWidth3 = $shape3.Width
$origHeight3 = $shape3.Height
$tr3 = $shape3.TextFrame.TextRange
$tr3.Runs(1, 1).Text = "- Having hard time creating pitchdeck"
$tr3.Runs(2, 1).Text = "- "
$shape3.Width = $origWidth3
$shape3.Height = $origHeight3

# --- Slide 4 ---------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shape4 = $s4.Shapes.Item(3)
Set-ShapeRunText $shape4 2 "PitchItup - an auto generated pitchdeck"
